$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.444.41'
Set-TextValue $ws.Range("E2") '  -0.13%  '
Set-TextValue $ws.Range("D3") '1.852.38'
Set-TextValue $ws.Range("E3") '  -0.04%  '
Set-TextValue $ws.Range("E4") '  +0.16%  '
Set-TextValue $ws.Range("D5") '240.87'
Set-TextValue $ws.Range("E5") '  +0.09%  '
Set-TextValue $ws.Range("D6") '0.6308'
Set-TextValue $ws.Range("E6") '  -0.23%  '
Set-TextValue $ws.Range("E7") '  +0.10%  '
Set-TextValue $ws.Range("E8") '  +1.42%  '
Set-TextValue $ws.Range("D9") '0.2944'
Set-TextValue $ws.Range("E9") '  -0.64%  '
Set-TextValue $ws.Range("D10") '24.60'
Set-TextValue $ws.Range("E10") '  -0.34%  '
Set-TextValue $ws.Range("E11") '  +0.63%  '
Set-TextValue $ws.Range("D12") '1.851.81'
Set-TextValue $ws.Range("E12") '  -0.30%  '
Set-TextValue $ws.Range("E13") '  +8.46%  '
Set-TextValue $ws.Range("D14") '5.028'
Set-TextValue $ws.Range("E14") '  +0.46%  '
Set-TextValue $ws.Range("D15") '0.6810'
Set-TextValue $ws.Range("E15") '  -0.48%  '
Set-TextValue $ws.Range("D16") '83.60'
Set-TextValue $ws.Range("E16") '  +0.26%  '
Set-TextValue $ws.Range("D17") '2.105.67'
Set-TextValue $ws.Range("E17") '  +0.43%  '
Set-TextValue $ws.Range("D18") '6.166'
Set-TextValue $ws.Range("E18") '  +0.17%  '
Set-TextValue $ws.Range("D19") '29.464.76'
Set-TextValue $ws.Range("E19") '  -0.13%  '
Set-TextValue $ws.Range("D20") '229.48'
Set-TextValue $ws.Range("E20") '  +0.20%  '
Set-TextValue $ws.Range("E21") '  -0.40%  '
Set-TextValue $ws.Range("E22") '  +0.11%  '
Set-TextValue $ws.Range("D23") '7.458'
Set-TextValue $ws.Range("E23") '  -0.96%  '
Set-TextValue $ws.Range("E24") '  +0.13%  '
Set-TextValue $ws.Range("D25") '157.33'
Set-TextValue $ws.Range("E25") '  +0.33%  '
Set-TextValue $ws.Range("E26") '  -0.98%  '
Set-TextValue $ws.Range("D27") '8.404'
Set-TextValue $ws.Range("E27") '  +0.16%  '
Set-TextValue $ws.Range("D28") '17.72'
Set-TextValue $ws.Range("E28") '  +0.20%  '
Set-TextValue $ws.Range("D29") '1.322'
Set-TextValue $ws.Range("E29") '  +3.85%  '
Set-TextValue $ws.Range("D30") '1.469'
Set-TextValue $ws.Range("E30") '  +0.10%  '
Set-TextValue $ws.Range("D31") '0.05697'
Set-TextValue $ws.Range("E31") '  +0.22%  '
Set-TextValue $ws.Range("E32") '  +0.07%  '
Set-TextValue $ws.Range("E33") '  +0.38%  '
Set-TextValue $ws.Range("D34") '1.853'
Set-TextValue $ws.Range("E34") '  +0.27%  '
Set-TextValue $ws.Range("D35") '1.162'
Set-TextValue $ws.Range("E35") '  +0.14%  '
Set-TextValue $ws.Range("D36") '0.7104'
Set-TextValue $ws.Range("E36") '  -1.89%  '
Set-TextValue $ws.Range("E37") '  -0.16%  '
Set-TextValue $ws.Range("D38") '2.783'
Set-TextValue $ws.Range("D39") '0.01795'
Set-TextValue $ws.Range("E39") '  -0.72%  '
Set-TextValue $ws.Range("D40") '1.217.91'
Set-TextValue $ws.Range("E40") '  -2.31%  '
Set-TextValue $ws.Range("E41") '  +4.83%  '
Set-TextValue $ws.Range("D42") '0.9078'
Set-TextValue $ws.Range("E42") '  -0.41%  '
Set-TextValue $ws.Range("E43") '  +0.10%  '
Set-TextValue $ws.Range("D44") '2.014.41'
Set-TextValue $ws.Range("E44") '  +0.42%  '
Set-TextValue $ws.Range("D45") '101.79'
Set-TextValue $ws.Range("E45") '  +0.00%  '
Set-TextValue $ws.Range("D46") '66.57'
Set-TextValue $ws.Range("E46") '  +0.25%  '
Set-TextValue $ws.Range("D47") '0.00000000119'
Set-TextValue $ws.Range("E47") '  -0.22%  '
Set-TextValue $ws.Range("D48") '7.142'
Set-TextValue $ws.Range("E48") '  +0.94%  '
Set-TextValue $ws.Range("D49") '0.4021'
Set-TextValue $ws.Range("E49") '  -0.22%  '
Set-TextValue $ws.Range("D50") '9.061'
Set-TextValue $ws.Range("E50") '  -0.55%  '
Set-TextValue $ws.Range("D51") '1.686'
Set-TextValue $ws.Range("E51") '  -0.81%  '
